# ----------------------------------------------------------------------
# NYPD CompStat weekly report refresh: roll the report forward one week
# (Volume/Number + date-range header) and update the precinct crime-stat
# grid (rows 15-30) with newly collected figures.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Volume 30, Number 33 -> 34 ---
$ws.Range("A8").Characters(21, 2).Text = "34"

# --- Header: report week 8/14/2023-8/20/2023 -> 8/21/2023-8/27/2023 ---
$ws.Range("C9").Characters(27, 9).Text = "8/21/2023"
$ws.Range("C9").Characters(47, 9).Text = "8/27/2023"

# --- Crime-stat grid updates (rows 15-30) ---
# Row 15
$ws.Range("J15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1
$ws.Range("F16").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("I15").Value = 9
$ws.Range("K15").Value = 80
$ws.Range("L15").Value = 28.571428571428
$ws.Range("M15").Value = 28.571428571428
$ws.Range("N15").Value = -59.090909090909
# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -26.666666666666
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = -31.034482758620
$ws.Range("L16").Value = -11.111111111111
$ws.Range("M16").Value = -42.028985507246
$ws.Range("N16").Value = -86.463620981387
# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -25
$ws.Range("I17").Value = 125
$ws.Range("J17").Value = 123
$ws.Range("K17").Value = 1.626016260162
$ws.Range("L17").Value = 17.924528301886
$ws.Range("M17").Value = 78.571428571428
$ws.Range("N17").Value = -62.235649546827
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -57.142857142857
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -43.75
$ws.Range("I18").Value = 119
$ws.Range("J18").Value = 123
$ws.Range("K18").Value = -3.252032520325
$ws.Range("L18").Value = 45.121951219512
$ws.Range("M18").Value = 43.373493975903
$ws.Range("N18").Value = -82.828282828282
# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -16.326530612244
$ws.Range("I19").Value = 331
$ws.Range("J19").Value = 340
$ws.Range("K19").Value = -2.647058823529
$ws.Range("L19").Value = 19.927536231884
$ws.Range("M19").Value = -2.359882005899
$ws.Range("N19").Value = -52.646638054363
# Row 20
$ws.Range("J20").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 2
$ws.Range("A15").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "***.*"
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 69
$ws.Range("K20").Value = 25.454545454545
$ws.Range("L20").Value = 18.965517241379
$ws.Range("M20").Value = 155.555555555556
$ws.Range("N20").Value = -89.497716894977
# Row 21
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = -21.904761904761
$ws.Range("I21").Value = 735
$ws.Range("J21").Value = 763
$ws.Range("K21").Value = -3.669724770642
$ws.Range("L21").Value = 18.548387096774
$ws.Range("M21").Value = 10.360360360360
$ws.Range("N21").Value = -75.581395348837
# Row 22
$ws.Range("E15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "0"
$ws.Range("F22").Value = 3
$ws.Range("G15").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value = "0"
$ws.Range("H15").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H22").Value = "***.*"
# Row 23
$ws.Range("C23").Value = 3
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 12
$ws.Range("H23").Value = 71.428571428571
$ws.Range("I23").Value = 77
$ws.Range("J23").Value = 84
$ws.Range("K23").Value = -8.333333333333
$ws.Range("L23").Value = 1.315789473684
$ws.Range("M23").Value = 45.283018867924
# Row 24
$ws.Range("C24").Value = 28
$ws.Range("E24").Value = -34.883720930232
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 168
$ws.Range("H24").Value = -45.833333333333
$ws.Range("I24").Value = 1080
$ws.Range("J24").Value = 1207
$ws.Range("K24").Value = -10.521955260977
$ws.Range("L24").Value = 54.727793696275
$ws.Range("M24").Value = 60.714285714285
# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 52.941176470588
$ws.Range("I25").Value = 197
$ws.Range("J25").Value = 199
$ws.Range("K25").Value = -1.005025125628
$ws.Range("L25").Value = 13.872832369942
$ws.Range("M25").Value = -8.796296296296
# Row 26
$ws.Range("I22").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = 1
$ws.Range("J22").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Value = 1
$ws.Range("I26").Value = 12
$ws.Range("K26").Value = 71.428571428571
$ws.Range("L26").Value = 20
# Row 27
$ws.Range("D23").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("G23").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("K22").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 32
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -3.030303030303
$ws.Range("L27").Value = -8.571428571428
# Row 28
$ws.Range("A16").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = "0"
$ws.Range("H28").Value = -100
$ws.Range("M28").Value = 150
$ws.Range("N28").Value = -85.294117647058
# Row 29
$ws.Range("A17").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = "0"
$ws.Range("H29").Value = -100
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -87.878787878787
# Row 30
$ws.Range("D24").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("L22").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("J26").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1
$ws.Range("M22").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 14
$ws.Range("K30").Value = -28.571428571428
